$wb = $excel.ActiveWorkbook

# --- Update the "tipo_movimiento" notes text (shared string) on tbakardez sheet ---
$ws1 = $wb.Worksheets.Item("tbakardez")
$ws1.Range("C2").Value = "E= Entrada, S= salida, A= Ajuste, V= Venta, T= Traspaso, C= Compra, I=Traspaso Ingreso,         D=Cancelar Venta"

# --- Widen column C slightly to accommodate the longer note text ---
$ws1.Columns.Item(3).ColumnWidth = 10.65

# --- Increase row 2 height so the wrapped note text fits ---
$ws1.Rows.Item(2).RowHeight = 165

# --- Move the active selection to C3 (as it ends up after editing C2) ---
$ws1.Range("C3").Select()

